# Robot Facturación Materia Prima / Config / Dia_Ejecucion.xlsx
# "Cambio transacción ZMM023, si no encuentra tickets de registro,
#  avanzar a procesar tickets en estado REV."
#
# - Refresh the execution date in A2 to the new run date.
# - Clear out the two stale historical date entries in A3/A4.
# - Leave the selection where the author left it (C9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the application window to match the saved workbook view
# (best effort / cosmetic - mirrors the author's window size at save time).
$excel.ActiveWindow.Width = 1536
$excel.ActiveWindow.Height = 606

# Update the execution date cell (A2) with the new date.
$ws.Range("A2").Value = "18/01/2023"

# Remove the old/no-longer-needed date entries.
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Match the final cell selection recorded in the workbook.
$ws.Range("C9").Select() | Out-Null
